# Refresh the cryptos price/volume snapshot (GitHub Actions data pull).
# Column D (Price) holds plain-text numbers (e.g. "2.00", "0.380",
# "61.091.22") whose exact formatting -- trailing zeros, dotted thousands --
# must be preserved, so each D-column cell is briefly switched to Text format
# before the value is written (otherwise Excel auto-coerces it to a number and
# drops the formatting), then the temporary format is cleared again so the
# cell keeps the workbook default style, matching the source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.091.22"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.641.97"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.58%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.81"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.19"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.69%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.599"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.55"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.108"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.20%  "
$ws.Range("E11").Value = "  +1.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.380"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.108.34"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.23"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +5.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "61.151.20"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000145"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.651.75"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.65"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.75"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "352.85"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.87"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.529"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.14"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.58"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +6.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.163"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.33%  "
$ws.Range("E27").Value = "  -0.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.00"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +6.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0818"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.89"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +7.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "168.39"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.05%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.03"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.85%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.10"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +10.74%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.66"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +7.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.36"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +8.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.70"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +4.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "341.21"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +9.90%  "
$ws.Range("B39").Value = "SuiNetwork"
$ws.Range("C39").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.927"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +10.34%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.14"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +5.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "38.22"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.32"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +5.91%  "
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0578"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +5.19%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "135.89"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.04"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +4.47%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "20.49"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.45%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.626"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0251"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +4.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0999"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.996"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.091.83"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.62%  "
